# Add 2022-Q1 data (commit: "feat: add 2022-Q1 data")
#
# 1. Insert a new worksheet "2022-Q1" right after "2021-Q4" (before "总计")
#    and fill it with the per-fund holdings table for that quarter.
# 2. Prepend a matching summary row to the "总计" (totals) sheet and
#    renumber its leading index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. New "2022-Q1" sheet
# ---------------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$ws = $wb.Worksheets.Add($null, $q4)
$ws.Name = "2022-Q1"

# Pick up the same header / index-column formatting (bold, centered, boxed)
# used by every other quarterly sheet.
$q4.Range("B1:H1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)
$q4.Range("A2").Copy()
$ws.Range("A2:A12").PasteSpecial(-4122)

$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

$rows = @(
    @("001955", "中欧养老产业混合",               "46.95", "92.35", "7.18", "3.3710", 8),
    @("010429", "中欧睿见混合",                   "27.61", "92.03", "7.27", "2.0072", 8),
    @("004616", "中欧电子信息产业沪港深股票A",       "14.54", "92.26", "4.75", "0.6906", 6),
    @("005763", "中欧电子信息产业沪港深股票C",       "7.73",  "92.26", "4.75", "0.3672", 6),
    @("003713", "英大睿盛灵活配置混合A",            "5.99",  "87.42", "4.14", "0.2480", 8),
    @("003714", "英大睿盛灵活配置混合C",            "2.40",  "87.42", "4.14", "0.0994", 8),
    @("014339", "长江智能制造混合A",               "3.28",  "21.63", "2.23", "0.0731", 2),
    @("003279", "融通沪港深智慧生活灵活配置混合",     "0.13",  "67.65", "3.62", "0.0047", 3),
    @("014340", "长江智能制造混合C",               "0.15",  "21.63", "2.23", "0.0033", 2),
    @("005146", "兴银丰润灵活配置混合",             "0.05",  "93.36", "3.17", "0.0016", 7),
    @("001608", "英大策略优选混合C",               "0.03",  "89.86", "5.15", "0.0015", 6)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $r - 2

    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 2).Value = $row[0]

    $ws.Cells.Item($r, 3).NumberFormat = "@"
    $ws.Cells.Item($r, 3).Value = $row[1]

    $ws.Cells.Item($r, 4).NumberFormat = "@"
    $ws.Cells.Item($r, 4).Value = $row[2]

    $ws.Cells.Item($r, 5).NumberFormat = "@"
    $ws.Cells.Item($r, 5).Value = $row[3]

    $ws.Cells.Item($r, 6).NumberFormat = "@"
    $ws.Cells.Item($r, 6).Value = $row[4]

    $ws.Cells.Item($r, 7).NumberFormat = "@"
    $ws.Cells.Item($r, 7).Value = $row[5]

    $ws.Cells.Item($r, 8).Value = $row[6]

    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 2. Update the "总计" (totals) sheet with the new quarter on top
# ---------------------------------------------------------------------------
$tot = $wb.Worksheets.Item("总计")

$tot.Rows.Item(2).Insert()

# The freshly inserted row doesn't inherit the bold/boxed index-column style
# (or the plain data style for B:D) - copy it over from the row below, which
# still carries the formatting every other data row in this sheet uses.
$tot.Range("A3:D3").Copy()
$tot.Range("A2:D2").PasteSpecial(-4122)

$tot.Cells.Item(2, 2).Value = "2022-Q1"
$tot.Cells.Item(2, 3).Value = 11
$tot.Cells.Item(2, 4).Value = 6.87

# Renumber the leading 0-based index column now that a row was inserted.
# (anchor on column B - xlDown from an empty A1 would not land correctly)
$last = $tot.Range("B1").End(-4121).Row
for ($i = 2; $i -le $last; $i++) {
    $tot.Cells.Item($i, 1).Value = $i - 2
}
